$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 463.732605
$ws.Range("H2").Value = 1391.197815
$ws.Range("I2").Value = 0.3632113435366598
$ws.Range("J2").Value = 0.3632113435366598
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.752878
$ws.Range("N2").Value = 5.258634
$ws.Range("O2").Value = 0.1377607590022273
$ws.Range("P2").Value = 0.1377607590022273
$ws.Range("Q2").Value = 812.8666811871899
$ws.Range("R2").Value = 7315.80013068471
$ws.Range("S2").Value = 0.05003627036382898
$ws.Range("T2").Value = 0.05003627036382898

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 463.732605
$ws.Range("H3").Value = 1391.197815
$ws.Range("I3").Value = 0.3632113435366598
$ws.Range("J3").Value = 0.3632113435366598
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.076282333333333
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.2417686736584878
$ws.Range("P3").Value = 0.2417686736584878
$ws.Range("Q3").Value = 1426.572420152145
$ws.Range("R3").Value = 12839.1517813693
$ws.Range("S3").Value = 0.0878131247845756
$ws.Range("T3").Value = 0.08781312478457562

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 463.732605
$ws.Range("H4").Value = 1391.197815
$ws.Range("I4").Value = 0.3632113435366598
$ws.Range("J4").Value = 0.3632113435366598
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.303088666666667
$ws.Range("N4").Value = 21.909266
$ws.Range("O4").Value = 0.5739583917309499
$ws.Range("P4").Value = 0.5739583917309499
$ws.Range("Q4").Value = 3386.68033193931
$ws.Range("R4").Value = 30480.12298745379
$ws.Range("S4").Value = 0.2084681985947388
$ws.Range("T4").Value = 0.2084681985947388

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 463.732605
$ws.Range("H5").Value = 1391.197815
$ws.Range("I5").Value = 0.3632113435366598
$ws.Range("J5").Value = 0.3632113435366598
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5918243333333334
$ws.Range("N5").Value = 1.775473
$ws.Range("O5").Value = 0.04651217560833507
$ws.Range("P5").Value = 0.04651217560833507
$ws.Range("Q5").Value = 274.448239799055
$ws.Range("R5").Value = 2470.034158191495
$ws.Range("S5").Value = 0.01689374979351644
$ws.Range("T5").Value = 0.01689374979351644

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.24662533333333
$ws.Range("H6").Value = 36.739876
$ws.Range("I6").Value = 0.009591978638444229
$ws.Range("J6").Value = 0.009591978638444227
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.752878
$ws.Range("N6").Value = 5.258634
$ws.Range("O6").Value = 0.1377607590022273
$ws.Range("P6").Value = 0.1377607590022273
$ws.Range("Q6").Value = 21.46684012104267
$ws.Range("R6").Value = 193.201561089384
$ws.Range("S6").Value = 0.001321398257565228
$ws.Range("T6").Value = 0.001321398257565227

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.24662533333333
$ws.Range("H7").Value = 36.739876
$ws.Range("I7").Value = 0.009591978638444229
$ws.Range("J7").Value = 0.009591978638444227
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.076282333333333
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.2417686736584878
$ws.Range("P7").Value = 0.2417686736584878
$ws.Range("Q7").Value = 37.67407715588578
$ws.Range("R7").Value = 339.066694402972
$ws.Range("S7").Value = 0.002319039953177209
$ws.Range("T7").Value = 0.002319039953177209

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.24662533333333
$ws.Range("H8").Value = 36.739876
$ws.Range("I8").Value = 0.009591978638444229
$ws.Range("J8").Value = 0.009591978638444227
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.303088666666667
$ws.Range("N8").Value = 21.909266
$ws.Range("O8").Value = 0.5739583917309499
$ws.Range("P8").Value = 0.5739583917309499
$ws.Range("Q8").Value = 89.43819067677956
$ws.Range("R8").Value = 804.943716091016
$ws.Range("S8").Value = 0.005505396632839076
$ws.Range("T8").Value = 0.005505396632839075

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.24662533333333
$ws.Range("H9").Value = 36.739876
$ws.Range("I9").Value = 0.009591978638444229
$ws.Range("J9").Value = 0.009591978638444227
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5918243333333334
$ws.Range("N9").Value = 1.775473
$ws.Range("O9").Value = 0.04651217560833507
$ws.Range("P9").Value = 0.04651217560833507
$ws.Range("Q9").Value = 7.247850873483112
$ws.Range("R9").Value = 65.230657861348
$ws.Range("S9").Value = 0.0004461437948627167
$ws.Range("T9").Value = 0.0004461437948627166

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 689.7685036666667
$ws.Range("H10").Value = 2069.305511
$ws.Range("I10").Value = 0.5402504422695089
$ws.Range("J10").Value = 0.5402504422695089
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.752878
$ws.Range("N10").Value = 5.258634
$ws.Range("O10").Value = 0.1377607590022273
$ws.Range("P10").Value = 0.1377607590022273
$ws.Range("Q10").Value = 1209.080035170219
$ws.Range("R10").Value = 10881.72031653197
$ws.Range("S10").Value = 0.07442531097833652
$ws.Range("T10").Value = 0.07442531097833652

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 689.7685036666667
$ws.Range("H11").Value = 2069.305511
$ws.Range("I11").Value = 0.5402504422695089
$ws.Range("J11").Value = 0.5402504422695089
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.076282333333333
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.2417686736584878
$ws.Range("P11").Value = 0.2417686736584878
$ws.Range("Q11").Value = 2121.922661919535
$ws.Range("R11").Value = 19097.30395727582
$ws.Range("S11").Value = 0.1306156328709106
$ws.Range("T11").Value = 0.1306156328709106

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 689.7685036666667
$ws.Range("H12").Value = 2069.305511
$ws.Range("I12").Value = 0.5402504422695089
$ws.Range("J12").Value = 0.5402504422695089
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 7.303088666666667
$ws.Range("N12").Value = 21.909266
$ws.Range("O12").Value = 0.5739583917309499
$ws.Range("P12").Value = 0.5739583917309499
$ws.Range("Q12").Value = 5037.440541751658
$ws.Range("R12").Value = 45336.96487576493
$ws.Range("S12").Value = 0.3100812749769417
$ws.Range("T12").Value = 0.3100812749769417

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 689.7685036666667
$ws.Range("H13").Value = 2069.305511
$ws.Range("I13").Value = 0.5402504422695089
$ws.Range("J13").Value = 0.5402504422695089
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5918243333333334
$ws.Range("N13").Value = 1.775473
$ws.Range("O13").Value = 0.04651217560833507
$ws.Range("P13").Value = 0.04651217560833507
$ws.Range("Q13").Value = 408.2217848368559
$ws.Range("R13").Value = 3673.996063531703
$ws.Range("S13").Value = 0.02512822344332009
$ws.Range("T13").Value = 0.02512822344332009

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 111.00921
$ws.Range("H14").Value = 333.02763
$ws.Range("I14").Value = 0.08694623555538696
$ws.Range("J14").Value = 0.08694623555538696
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.752878
$ws.Range("N14").Value = 5.258634
$ws.Range("O14").Value = 0.1377607590022273
$ws.Range("P14").Value = 0.1377607590022273
$ws.Range("Q14").Value = 194.58560200638
$ws.Range("R14").Value = 1751.27041805742
$ws.Range("S14").Value = 0.01197777940249655
$ws.Range("T14").Value = 0.01197777940249655

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 111.00921
$ws.Range("H15").Value = 333.02763
$ws.Range("I15").Value = 0.08694623555538696
$ws.Range("J15").Value = 0.08694623555538696
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.076282333333333
$ws.Range("N15").Value = 9.228847
$ws.Range("O15").Value = 0.2417686736584878
$ws.Range("P15").Value = 0.2417686736584878
$ws.Range("Q15").Value = 341.49567156029
$ws.Range("R15").Value = 3073.46104404261
$ws.Range("S15").Value = 0.02102087604982436
$ws.Range("T15").Value = 0.02102087604982436

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 111.00921
$ws.Range("H16").Value = 333.02763
$ws.Range("I16").Value = 0.08694623555538696
$ws.Range("J16").Value = 0.08694623555538696
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 7.303088666666667
$ws.Range("N16").Value = 21.909266
$ws.Range("O16").Value = 0.5739583917309499
$ws.Range("P16").Value = 0.5739583917309499
$ws.Range("Q16").Value = 810.71010344662
$ws.Range("R16").Value = 7296.39093101958
$ws.Range("S16").Value = 0.04990352152643023
$ws.Range("T16").Value = 0.04990352152643023

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 111.00921
$ws.Range("H17").Value = 333.02763
$ws.Range("I17").Value = 0.08694623555538696
$ws.Range("J17").Value = 0.08694623555538696
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5918243333333334
$ws.Range("N17").Value = 1.775473
$ws.Range("O17").Value = 0.04651217560833507
$ws.Range("P17").Value = 0.04651217560833507
$ws.Range("Q17").Value = 65.69795170211
$ws.Range("R17").Value = 591.28156531899
$ws.Range("S17").Value = 0.004044058576635825
$ws.Range("T17").Value = 0.004044058576635825

